# V 0.52-B49 - Add Ambient Vertical Wind Item (LIFT)
# Inserts a new "VWIND" data column into Tabelle2 right before the
# END_OF_COL / Title columns, shifting the trailing columns one to the
# right, and registers the new shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle2")

# Column 134 == "ED" -- this is where END_OF_COL currently lives.
# Inserting here pushes END_OF_COL (-> EE) and Title (-> EF) right by one
# column and creates a fresh column ED for the new VWIND data.
$ws.Columns.Item(134).Insert()

# Header cell for the new column.
$ws.Cells.Item(1, 134).Value = "VWIND"

# Data rows 2-40 get the same "|" placeholder used by the other flag
# style columns in this sheet.
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 134).Value = "|"
}

# Restore the selected cell shown in the saved view.
$ws.Range("EL3").Select()
